$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text storage for numeric-looking price values so they remain strings (matching source data)
$priceCells = @("D2", "D4", "D5", "D7", "D8", "D9", "D10", "D11", "D12", "D13", "D14", "D15", "D16", "D17", "D18", "D19", "D20", "D21", "D23", "D24", "D27", "D40", "D41", "D42", "D43", "D44", "D45", "D47", "D48")
foreach ($c in $priceCells) { $ws.Range($c).NumberFormat = "@" }

# Apply the updated cell values
$ws.Range("D2").Value = "246.32"

$ws.Range("D4").Value = "5.423"

$ws.Range("D5").Value = "0.05780"

$ws.Range("D7").Value = "6.326"

$ws.Range("D8").Value = "0.8171"

$ws.Range("D9").Value = "0.9467"
$ws.Range("E9").Value = "8FTXTokenFTTBestin24h"

$ws.Range("B10").Value = "WazirX"
$ws.Range("C10").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("D10").Value = "0.1428"
$ws.Range("E10").Value = "9WazirXWRX"

$ws.Range("B11").Value = "MandalaExchangeToken"
$ws.Range("C11").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("D11").Value = "0.07475"
$ws.Range("E11").Value = "10MandalaExchangeTokenMDX"

$ws.Range("B12").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C12").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$ws.Range("D12").Value = "0.03140"
$ws.Range("E12").Value = "11LiechtensteinCryptoassetsExchangeLCX"

$ws.Range("B13").Value = "BitrueCoin"
$ws.Range("C13").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("D13").Value = "0.03001"
$ws.Range("E13").Value = "12BitrueCoinBTR"

$ws.Range("B14").Value = "MCDex"
$ws.Range("C14").Value = "https://coinranking.com/coin/3nMM61qeg+mcdex-mcb"
$ws.Range("D14").Value = "4.143"
$ws.Range("E14").Value = "13MCDexMCB"

$ws.Range("B15").Value = "BitMartToken"
$ws.Range("C15").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("D15").Value = "0.09408"
$ws.Range("E15").Value = "14BitMartTokenBMX"

$ws.Range("B16").Value = "BitForexToken"
$ws.Range("C16").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("D16").Value = "0.001591"
$ws.Range("E16").Value = "15BitForexTokenBF"

$ws.Range("B17").Value = "CoinExToken"
$ws.Range("C17").Value = "https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet"
$ws.Range("D17").Value = "0.04829"
$ws.Range("E17").Value = "16CoinExTokenCET"

$ws.Range("B18").Value = "One"
$ws.Range("C18").Value = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
$ws.Range("D18").Value = "0.0005851"
$ws.Range("E18").Value = "17OneONE"

$ws.Range("D19").Value = "0.006196"

$ws.Range("D20").Value = "0.004123"

$ws.Range("D21").Value = "0.0009974"

$ws.Range("D23").Value = "3.771"

$ws.Range("D24").Value = "2.223"

$ws.Range("D27").Value = "0.0003999"

$ws.Range("D40").Value = "0.03895"

$ws.Range("B41").Value = "BKEXToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
$ws.Range("D41").Value = "0.1077"
$ws.Range("E41").Value = "40BKEXTokenBKK"

$ws.Range("B42").Value = "CEJI"
$ws.Range("C42").Value = "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"
$ws.Range("D42").Value = "0.002621"
$ws.Range("E42").Value = "41CEJICEJI"

$ws.Range("B43").Value = "KickToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick"
$ws.Range("D43").Value = "0.003038"
$ws.Range("E43").Value = "42KickTokenKICK"

$ws.Range("D44").Value = "0.006548"

$ws.Range("D45").Value = "0.00005593"

$ws.Range("D47").Value = "0.3800"

$ws.Range("D48").Value = "0.1486"

# Restore default style on the price cells (NumberFormat change above would otherwise leave a stray style)
foreach ($c in $priceCells) { $ws.Range($c).Style = "Normal" }

Write-Host "Applied cryptos.xlsx symbol-list update"
